# Update OIEC surface analysis results worksheet: add a new "2PC BM analysis" block
# (rows 36-48) mirroring the two existing blocks above it, and retarget the row-33
# "theta" row height to match its twin in the new block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 33: height only changes (16.5 -> 14.7); content is untouched.
# ---------------------------------------------------------------------------
$ws.Rows(33).RowHeight = 14.7

# ---------------------------------------------------------------------------
# Row 36: new section title "2PC BM analysis" (clear any leftover formatting
# and custom height from the previously-blank placeholder row first).
# ---------------------------------------------------------------------------
$ws.Range("A36").Style = "Normal"
$ws.Range("A36").Value = "2PC BM analysis"
$ws.Rows(36).AutoFit()

# ---------------------------------------------------------------------------
# Row 37: "Parameter" / "Value" header, bold -- copy format from the header of
# the block above (row 22) which already carries the bold style.
# ---------------------------------------------------------------------------
$ws.Range("A37").Style = "Normal"
$ws.Range("A22:B22").Copy($ws.Range("A37:B37"))
$ws.Range("A37").Value = "Parameter"
$ws.Range("B37").Value = "Value"
$ws.Rows(37).AutoFit()

# ---------------------------------------------------------------------------
# Row 38: AICc value (plain, unstyled cells).
# ---------------------------------------------------------------------------
$ws.Range("A38").Value = "AICc"
$ws.Range("B38").Value = -173.88333800000001

# ---------------------------------------------------------------------------
# Rows 39-40: Phenotypic regimes / Pheno reg shifts, value "-" centered AND
# vertically centered (new style).
# ---------------------------------------------------------------------------
$ws.Range("A39").Value = "Phenotypic regimes"
$ws.Range("B39").Value = "-"
$ws.Range("B39").HorizontalAlignment = -4108
$ws.Range("B39").VerticalAlignment = -4108

$ws.Range("A40").Value = "Pheno reg shifts"
$ws.Range("B40").Value = "-"
$ws.Range("B40").HorizontalAlignment = -4108
$ws.Range("B40").VerticalAlignment = -4108

# ---------------------------------------------------------------------------
# Rows 41-43: Conv pheno reg / Conv reg shifts / Conv fraction, value "-"
# centered horizontally only (re-uses the pre-existing centered style).
# ---------------------------------------------------------------------------
$ws.Range("A41").Value = "Conv pheno reg"
$ws.Range("B41").Value = "-"
$ws.Range("B41").HorizontalAlignment = -4108

$ws.Range("A42").Value = "Conv reg shifts"
$ws.Range("B42").Value = "-"
$ws.Range("B42").HorizontalAlignment = -4108

$ws.Range("A43").Value = "Conv fraction"
$ws.Range("B43").Value = "-"
$ws.Range("B43").HorizontalAlignment = -4108

# ---------------------------------------------------------------------------
# Row 44: PC1 / PC2 column headers with the bottom-border style, copied from
# the equivalent header row of the block above (row 29).
# ---------------------------------------------------------------------------
$ws.Range("A29:C29").Copy($ws.Range("A44:C44"))

# ---------------------------------------------------------------------------
# Rows 45-48: parameter rows (alpha, t1/2, sigma^2, theta). Copy the label
# cell formatting (and, for rich-text labels, the exact shared string) from
# the matching row of the block above, then overwrite values as needed.
# ---------------------------------------------------------------------------
$ws.Range("A30").Copy($ws.Range("A45"))
$ws.Range("B45").Value = "-"
$ws.Range("B45").HorizontalAlignment = -4108
$ws.Range("C45").Value = "-"
$ws.Range("C45").HorizontalAlignment = -4108
$ws.Rows(45).RowHeight = 14.7

$ws.Range("A31").Copy($ws.Range("A46"))
$ws.Range("B46").Value = "-"
$ws.Range("B46").HorizontalAlignment = -4108
$ws.Range("C46").Value = "-"
$ws.Range("C46").HorizontalAlignment = -4108
$ws.Rows(46).RowHeight = 16.8

$ws.Range("A32:C32").Copy($ws.Range("A47:C47"))
$ws.Range("B47").Value = 0.00013899999999999999
$ws.Range("C47").Value = 0.00032000000000000003
$ws.Rows(47).RowHeight = 16.5

$ws.Range("A33").Copy($ws.Range("A48"))
$ws.Range("B48").Value = "-"
$ws.Range("B48").HorizontalAlignment = -4108
$ws.Range("C48").Value = "-"
$ws.Range("C48").HorizontalAlignment = -4108
$ws.Rows(48).RowHeight = 14.7

# ---------------------------------------------------------------------------
# Update the view: scroll so row 29 is at the top and select D51, matching
# where the author ended up after entering the new data.
# ---------------------------------------------------------------------------
$ws.Range("D51").Select()
$excel.ActiveWindow.ScrollRow = 29
$excel.ActiveWindow.ScrollColumn = 1
